$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply a custom number format (0.0000) to the Longitude/Latitude decimal columns (D2:E8)
$ws.Range("D2:E8").NumberFormat = "0.0000"

# Add the missing differential elevation data point for row 8 (Junction T4)
$ws.Range("D8").Value = 2.31968889
$ws.Range("E8").Value = 33.24888889

# Update the active selection to F1, matching the post-edit workbook state
$ws.Range("F1").Select()
